$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 6 for the new "must have 4 modes..." TODO item.
# This shifts the existing rows 6-14 down to 7-15.
$ws.Rows.Item(6).Insert()

# Set the new cell values (rows 16 and 17 are brand new, past the previous
# last row, so simply assigning to them appends the new rows).
# The order in which .Value is assigned controls the order new entries are
# appended to the shared-strings table, so we set the text that must land at
# shared-string index 12 first, then 13, then 14.
$ws.Range("A17").Value = "Auto-orient by EXIF + reset exif orient tag (make sure other metadata is left intact)"
$ws.Range("B6").Value = "must have 4 modes = one is ""do not modify - just renaming and adding to output mix"""
$ws.Range("A16").Value = "Add credits for CCR-Exif and NativeJpg"

# Apply cell styles matching the updated styles.xml (new "Good" and
# "Explanatory Text" cell styles). "Good" must be registered before
# "Explanatory Text" so the new style indexes line up (2 = Good, 3 =
# Explanatory Text).
$ws.Range("A12").Style = "Good"
$ws.Range("A14").Style = "Good"
$ws.Range("A17").Style = "Good"
$ws.Range("A8").Style = "Explanatory Text"

# Update the active selection to match the committed worksheet state.
$ws.Range("A16").Select()
